$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged) - update values
$ws.Range("B3").Value = 4368708010079.297
$ws.Range("C3").Value = 4257358309887.731
$ws.Range("D3").Value = 342974371019027.9

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 3903242085595.563
$ws.Range("C4").Value = 3829323594115.925
$ws.Range("D4").Value = 179997264516186.5

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 705232929002432
$ws.Range("C5").Value = 937431475126992.2
$ws.Range("D5").Value = 3405232554532662
